# The sharedStrings table was reordered/extended by the author (Result and
# SamplePortion swapped slots, and the "#float" type string was split into
# unit-qualified variants) while the sheet's cell references were left
# pointing at the same slot indexes - net effect: the H1/I1 header text
# swaps, H2/I2/M2 get unit-qualified type strings, and a new row 3 of
# French "#..." enum/description labels is appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - headers: swap SamplePortion/Result between columns H and I
$ws.Range("H1").Value = "SamplePortion"
$ws.Range("I1").Value = "Result"

# Row 2 - types: qualify the float types with their units
$ws.Range("H2").Value = "#float,  unit:mg"
$ws.Range("I2").Value = "#float,  unit:mg"
$ws.Range("M2").Value = "#float,  unit:celsius"

# Row 3 - new French description/enum row
$ws.Range("A3").Value = "#Manipulateur"
$ws.Range("B3").Value = "#Desc:IdentifiantEchantillon"
$ws.Range("C3").Value = "#Date"
$ws.Range("D3").Value = "#ModeOderatoireLaboratoire"
$ws.Range("E3").Value = "#AppareilLogicielCritique"
$ws.Range("F3").Value = "#ProduitCritique"
$ws.Range("G3").Value = "#LieuStockageDonneesBrutes"
$ws.Range("H3").Value = "#PriseEssai"
$ws.Range("I3").Value = "#Resultat"
$ws.Range("J3").Value = "#Solvant"
$ws.Range("K3").Value = "#CycleDeTemperature"
$ws.Range("L3").Value = "#TypeDeGodet"
$ws.Range("M3").Value = "#TemperatureDeChauffage"
